$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block of results for a heuristic method for finding initial variables
$ws.Range("A8").Value = "E-n101-k14"
$ws.Range("B8").Value = 1071
$ws.Range("C8").Value = 114
$ws.Range("D8").Value = "1023,2…"
$ws.Range("E8").Value = "278s"

$ws.Range("A10").Value = "X-n101-k25"
$ws.Range("B10").Value = 27591
$ws.Range("C10").Value = 140
$ws.Range("D10").Value = "26787,8…"
$ws.Range("E10").Value = "97s"

$ws.Range("A11").Value = "X-n129-k18"
$ws.Range("B11").Value = 28940
$ws.Range("C11").Value = 41

$ws.Range("D11").Select()
